$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new headers I1 ("I0") and J1 ("IF"), matching the style of the existing
# header row (bold, bordered, centered) by copying H1's format first.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1:J1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Fill in the new I0 / IF numeric columns for each data row.
$values = @{
    2 = @(3, 4)
    3 = @(10, 10)
    4 = @(7, 7)
    5 = @(9, 9)
    6 = @(8, 9)
    7 = @(6, 7)
    8 = @(6, 7)
    9 = @(10, 10)
    10 = @(7, 8)
    11 = @(6, 7)
    12 = @(6, 6)
    13 = @(9, 9)
    14 = @(9, 9)
    15 = @(7, 7)
    16 = @(8, 8)
    17 = @(7, 7)
    18 = @(7, 7)
    19 = @(6, 6)
    20 = @(7, 8)
    21 = @(9, 9)
    22 = @(7, 7)
    23 = @(6, 6)
    24 = @(9, 9)
    25 = @(7, 7)
    26 = @(9, 9)
    27 = @(6, 6)
    28 = @(10, 11)
}

foreach ($row in $values.Keys) {
    $pair = $values[$row]
    $ws.Cells.Item($row, 9).Value = $pair[0]
    $ws.Cells.Item($row, 10).Value = $pair[1]
}

Write-Host "I0/IF columns added"
